$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column D width: 12 -> 11 (stored OOXML width = ColumnWidth + 0.8333333333333334)
$ws.Columns.Item(4).ColumnWidth = 10.166666666666666

# Row 2 (OTROS)
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

# Row 3 (PORCELANATO)
$ws.Range("C3").Value = 20000
$ws.Range("E3").Value = 20000

# Row 4 (TOTAL)
$ws.Range("C4").Value = 20000
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 20000
$ws.Range("F4").Value = 0
